# Update "想去人数" (column F) counts across all sheets to match the
# regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetChanges = @{
    "展览" = @{
        2 = 769
        3 = 14332
        4 = 14466
        6 = 1409
        7 = 5924
        9 = 579
        13 = 1565
        15 = 2133
        16 = 1217
        17 = 1860
        18 = 919
        20 = 2295
        21 = 572
        22 = 824
        23 = 3371
        25 = 314
        26 = 2435
        27 = 606
        30 = 1810
        31 = 1081
        32 = 1424
        33 = 110
        34 = 153
        35 = 4932
        36 = 4907
        39 = 682
        40 = 694
        41 = 3310
        45 = 114
        46 = 95
        47 = 4445
        48 = 615
        49 = 299
    };
    "演出" = @{
        3 = 47
        4 = 123
        17 = 18
        19 = 91
        24 = 2
    };
    "本地生活" = @{
        2 = 7650
        3 = 251
        4 = 855
    };
    "全部类型" = @{
        2 = 769
        3 = 251
        4 = 855
        6 = 14332
        7 = 47
        8 = 1409
        9 = 5924
        11 = 123
        14 = 1565
        16 = 1217
        17 = 1860
        18 = 919
        20 = 572
        21 = 3371
        22 = 314
        23 = 606
        25 = 1810
        28 = 1424
        29 = 91
        30 = 110
        31 = 153
        32 = 4933
        33 = 4908
        37 = 682
        38 = 694
        39 = 3310
        43 = 2
        44 = 95
        45 = 4445
        46 = 615
        47 = 299
    }
}

foreach ($sheetName in $sheetChanges.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $sheetChanges[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}
